# Update column G ("K") values for rows 2-10 to reflect the regenerated
# save data (K instead of Strike#, regen std/mean, calc and write s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 3
    5  = 0
    6  = 0
    8  = 1
    9  = 1
    10 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
